# Applies the BOQ/bill-summary update described in the commit diff.
# Rows 8, 10, 11, 12 get new line-item data; rows 14/16 totals are updated
# to match. Cells that hold "numbers stored as text" (D, G, H columns in the
# item rows, and the total rows) are forced to Text format first so Excel
# keeps them as text instead of auto-converting to numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- Row 8: "Each"/switch line -> "P. point"/medium point line ---
$ws.Range("A8").Value = 'P. point'
$ws.Range("C8").Value = 15
Set-TextValue $ws.Range("D8") "3"
$ws.Range("E8").Value = 'Medium point (up to 6 mtr.)'
$ws.Range("F8").Value = 472
Set-TextValue $ws.Range("G8") "7080.00"

# --- Row 9: quantity-only update ---
$ws.Range("C9").Value = 60

# --- Row 10: "Mtr."/wiring line -> blank unit/LED batten fixture line ---
$ws.Range("A10").Value = ""
$ws.Range("C10").Value = 36
Set-TextValue $ws.Range("D10") "16.0"
$ws.Range("E10").Value = 'Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F10").Value = 0
Set-TextValue $ws.Range("G10") "0.00"

# --- Row 11: "Set"/earthing line -> "Each"/LED batten spec line ---
$ws.Range("A11").Value = 'Each'
$ws.Range("C11").Value = 39
Set-TextValue $ws.Range("D11") "27"
$ws.Range("E11").Value = '1170mm(+/-10%) LED batten with min. lumen output 2200 lm'
$ws.Range("F11").Value = 492
Set-TextValue $ws.Range("G11") "19188.00"

# --- Row 12: "Grand Total" placeholder line -> Distribution board line ---
$ws.Range("C12").Value = 38
Set-TextValue $ws.Range("D12") "18.0"
$ws.Range("E12").Value = 'Providing & Fixing of Recessed/surface mounting heavy duty horizontal type Double Door ( Metal / Glazed )Distribution board with Metal end box made out from Galvanized steel / CRCA sheet not less then 1.2 mm thick  conforming to IS-8623-1 & 3 /  IEC 61439- 1 & 3, powder painted complete with reversible door (for double door DB only )100 amp.  insulated copper bus bar/shorting link , copper neutral link, copper earth link , color coded interconnecting wire set  of suitable rating and din bar,masking sheet,  making internal DB  terminations with copper lugs, Ferrules,  detachable gland plate, including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'

# --- Row 14 / 16: grand totals reflecting the new line items ---
Set-TextValue $ws.Range("G14") "26268.00"
Set-TextValue $ws.Range("H14") "26268.00"
Set-TextValue $ws.Range("G16") "26268.00"
Set-TextValue $ws.Range("H16") "26268.00"

